$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet (B08K42K3DM -> B08K4353V1)
$ws.Name = "B08K4353V1"

# 2. Replace all 100 keyword values in column A with the new keyword list
$ws.Cells.Item(1, 1).Value = "neck compression wrap"
$ws.Cells.Item(2, 1).Value = "training bra sets"
$ws.Cells.Item(3, 1).Value = "wide waistband yoga"
$ws.Cells.Item(4, 1).Value = "girls seamless underwear"
$ws.Cells.Item(5, 1).Value = "sport bra"
$ws.Cells.Item(6, 1).Value = "seamless bra"
$ws.Cells.Item(7, 1).Value = "yoga legging"
$ws.Cells.Item(8, 1).Value = "workout waist"
$ws.Cells.Item(9, 1).Value = "yoga workout"
$ws.Cells.Item(10, 1).Value = "yoga bra"
$ws.Cells.Item(11, 1).Value = "sports bra for women"
$ws.Cells.Item(12, 1).Value = "sport apparel"
$ws.Cells.Item(13, 1).Value = "seamless bra set"
$ws.Cells.Item(14, 1).Value = "workout apparel"
$ws.Cells.Item(15, 1).Value = "workout bra"
$ws.Cells.Item(16, 1).Value = "workout outfit set"
$ws.Cells.Item(17, 1).Value = "workout set"
$ws.Cells.Item(18, 1).Value = "workout leggings for women high waist"
$ws.Cells.Item(19, 1).Value = "workout legging"
$ws.Cells.Item(20, 1).Value = "legging bra set"
$ws.Cells.Item(21, 1).Value = "workout outfit"
$ws.Cells.Item(22, 1).Value = "sports leggings women"
$ws.Cells.Item(23, 1).Value = "black legging"
$ws.Cells.Item(24, 1).Value = "sport bra set"
$ws.Cells.Item(25, 1).Value = "yoga sets for women"
$ws.Cells.Item(26, 1).Value = "high waist workout"
$ws.Cells.Item(27, 1).Value = "yoga set purple"
$ws.Cells.Item(28, 1).Value = "yoga bras for women"
$ws.Cells.Item(29, 1).Value = "2 pcs"
$ws.Cells.Item(30, 1).Value = "green 2"
$ws.Cells.Item(31, 1).Value = "pcs set"
$ws.Cells.Item(32, 1).Value = "bra set"
$ws.Cells.Item(33, 1).Value = "waist yoga"
$ws.Cells.Item(34, 1).Value = "outfit set"
$ws.Cells.Item(35, 1).Value = "seamless sports bra"
$ws.Cells.Item(36, 1).Value = "seamless workout leggings"
$ws.Cells.Item(37, 1).Value = "seamless yoga bra"
$ws.Cells.Item(38, 1).Value = "seamless yoga"
$ws.Cells.Item(39, 1).Value = "sports for women"
$ws.Cells.Item(40, 1).Value = "high waist yoga"
$ws.Cells.Item(41, 1).Value = "sport legging"
$ws.Cells.Item(42, 1).Value = "black bra"
$ws.Cells.Item(43, 1).Value = "yoga legging set"
$ws.Cells.Item(44, 1).Value = "yoga workout set"
$ws.Cells.Item(45, 1).Value = "yoga bra set"
$ws.Cells.Item(46, 1).Value = "yoga outfit set"
$ws.Cells.Item(47, 1).Value = "workout set women"
$ws.Cells.Item(48, 1).Value = "black legging set"
$ws.Cells.Item(49, 1).Value = "sports leggings for women"
$ws.Cells.Item(50, 1).Value = "black bra set"
$ws.Cells.Item(51, 1).Value = "sport leggings for women"
$ws.Cells.Item(52, 1).Value = "sport set"
$ws.Cells.Item(53, 1).Value = "black outfit"
$ws.Cells.Item(54, 1).Value = "high waist"
$ws.Cells.Item(55, 1).Value = "set 2"
$ws.Cells.Item(56, 1).Value = "yoga apparel"
$ws.Cells.Item(57, 1).Value = "purple bra"
$ws.Cells.Item(58, 1).Value = "purple legging"
$ws.Cells.Item(59, 1).Value = "high waist legging"
$ws.Cells.Item(60, 1).Value = "high sport bra"
$ws.Cells.Item(61, 1).Value = "black n"
$ws.Cells.Item(62, 1).Value = "womens sports bra set"
$ws.Cells.Item(63, 1).Value = "sport outfit"
$ws.Cells.Item(64, 1).Value = "sport waist"
$ws.Cells.Item(65, 1).Value = "sport workout"
$ws.Cells.Item(66, 1).Value = "waist bra"
$ws.Cells.Item(67, 1).Value = "black 2"
$ws.Cells.Item(68, 1).Value = "yoga sport"
$ws.Cells.Item(69, 1).Value = "n set"
$ws.Cells.Item(70, 1).Value = "green set"
$ws.Cells.Item(71, 1).Value = "yoga 2"
$ws.Cells.Item(72, 1).Value = "green bra"
$ws.Cells.Item(73, 1).Value = "green bra set"
$ws.Cells.Item(74, 1).Value = "green yoga"
$ws.Cells.Item(75, 1).Value = "yoga sport bra"
$ws.Cells.Item(76, 1).Value = "womens workout sports bra"
$ws.Cells.Item(77, 1).Value = "high waist set"
$ws.Cells.Item(78, 1).Value = "ready set"
$ws.Cells.Item(79, 1).Value = "high black"
$ws.Cells.Item(80, 1).Value = "waist set"
$ws.Cells.Item(81, 1).Value = "black yoga"
$ws.Cells.Item(82, 1).Value = "green outfit"
$ws.Cells.Item(83, 1).Value = "black apparel"
$ws.Cells.Item(84, 1).Value = "purple yoga"
$ws.Cells.Item(85, 1).Value = "black set"
$ws.Cells.Item(86, 1).Value = "purple outfit"
$ws.Cells.Item(87, 1).Value = "purple workout"
$ws.Cells.Item(88, 1).Value = "black sport"
$ws.Cells.Item(89, 1).Value = "seamless legging"
$ws.Cells.Item(90, 1).Value = "green apparel"
$ws.Cells.Item(91, 1).Value = "purple set"
$ws.Cells.Item(92, 1).Value = "black seamless bra"
$ws.Cells.Item(93, 1).Value = "black yoga bra"
$ws.Cells.Item(94, 1).Value = "seamless workout"
$ws.Cells.Item(95, 1).Value = "green sport bra"
$ws.Cells.Item(96, 1).Value = "black sport bra"
$ws.Cells.Item(97, 1).Value = "purple sport bra"
$ws.Cells.Item(98, 1).Value = "workout sport bra"
$ws.Cells.Item(99, 1).Value = "black workout bra"
$ws.Cells.Item(100, 1).Value = "seamless sport bra"

# 3. Remove explicit "left" horizontal alignment from the data range (revert to default/general)
$ws.Range("A1:A100").HorizontalAlignment = 1

# 4. Re-assert the "duplicate values" conditional-formatting highlight on A1:A100.
#    Re-applying it (add + recolor + delete) twice mirrors the authors edit, which left two
#    extra (now-unused) duplicate fill entries in the style table alongside the original rule.
for ($i = 0; $i -lt 2; $i++) {
    $fc = $ws.Range("A1:A100").FormatConditions.AddUniqueValues()
    $fc.DupeUnique = 1
    $fc.Interior.Color = 13998939
    $fc.Delete()
}

# 5. Update the view state: scroll so row 10 is at the top, and select the full data range
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$ws.Range("A1:A100").Select()
